$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.185.33'
Set-TextValue $ws.Range('E2') '  +0.12%  '
Set-TextValue $ws.Range('D3') '1.838.65'
Set-TextValue $ws.Range('E3') '  +0.10%  '
Set-TextValue $ws.Range('D4') '1.009'
Set-TextValue $ws.Range('E4') '  +0.65%  '
Set-TextValue $ws.Range('D5') '243.73'
Set-TextValue $ws.Range('E5') '  -0.37%  '
Set-TextValue $ws.Range('D6') '0.6177'
Set-TextValue $ws.Range('E6') '  -1.93%  '
Set-TextValue $ws.Range('E7') '  +0.48%  '
Set-TextValue $ws.Range('D8') '0.07438'
Set-TextValue $ws.Range('E8') '  -0.88%  '
Set-TextValue $ws.Range('D9') '0.2942'
Set-TextValue $ws.Range('E9') '  +0.31%  '
Set-TextValue $ws.Range('D10') '22.93'
Set-TextValue $ws.Range('E10') '  -0.87%  '
Set-TextValue $ws.Range('D11') '0.07713'
Set-TextValue $ws.Range('E11') '  -0.41%  '
Set-TextValue $ws.Range('D12') '1.828.03'
Set-TextValue $ws.Range('E12') '  -0.51%  '
Set-TextValue $ws.Range('D13') '4.980'
Set-TextValue $ws.Range('E13') '  -0.23%  '
Set-TextValue $ws.Range('D14') '0.6709'
Set-TextValue $ws.Range('E14') '  +0.15%  '
Set-TextValue $ws.Range('D15') '82.73'
Set-TextValue $ws.Range('E15') '  -0.25%  '
Set-TextValue $ws.Range('D16') '0.000009082'
Set-TextValue $ws.Range('E16') '  -2.43%  '
Set-TextValue $ws.Range('D17') '5.870'
Set-TextValue $ws.Range('E17') '  -2.86%  '
Set-TextValue $ws.Range('D18') '29.203.88'
Set-TextValue $ws.Range('E18') '  +0.06%  '
Set-TextValue $ws.Range('D19') '2.086.11'
Set-TextValue $ws.Range('E19') '  -0.05%  '
Set-TextValue $ws.Range('D20') '236.62'
Set-TextValue $ws.Range('E20') '  +5.76%  '
Set-TextValue $ws.Range('D21') '12.59'
Set-TextValue $ws.Range('E21') '  -0.32%  '
Set-TextValue $ws.Range('E22') '  +0.43%  '
Set-TextValue $ws.Range('D23') '7.154'
Set-TextValue $ws.Range('D24') '1.013'
Set-TextValue $ws.Range('E24') '  +0.86%  '
Set-TextValue $ws.Range('D25') '159.50'
Set-TextValue $ws.Range('E25') '  -0.86%  '
Set-TextValue $ws.Range('D26') '0.1417'
Set-TextValue $ws.Range('E26') '  +1.25%  '
Set-TextValue $ws.Range('D27') '8.499'
Set-TextValue $ws.Range('E27') '  -0.32%  '
Set-TextValue $ws.Range('D28') '17.85'
Set-TextValue $ws.Range('E28') '  -0.66%  '
Set-TextValue $ws.Range('D29') '1.504'
Set-TextValue $ws.Range('E29') '  -0.11%  '
Set-TextValue $ws.Range('D30') '4.130'
Set-TextValue $ws.Range('E30') '  -0.71%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D31') '4.111'
Set-TextValue $ws.Range('E31') '  +0.83%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D32') '0.05553'
Set-TextValue $ws.Range('E32') '  -0.93%  '
Set-TextValue $ws.Range('D33') '1.217'
Set-TextValue $ws.Range('E33') '  +0.54%  '
Set-TextValue $ws.Range('D34') '1.853'
Set-TextValue $ws.Range('E34') '  -0.20%  '
Set-TextValue $ws.Range('D35') '0.7402'
Set-TextValue $ws.Range('E35') '  -1.64%  '
Set-TextValue $ws.Range('D36') '1.139'
Set-TextValue $ws.Range('E36') '  -0.05%  '
Set-TextValue $ws.Range('D37') '2.657'
Set-TextValue $ws.Range('E37') '  +1.50%  '
Set-TextValue $ws.Range('D38') '2.829'
Set-TextValue $ws.Range('E38') '  +2.77%  '
Set-TextValue $ws.Range('D39') '0.01773'
Set-TextValue $ws.Range('E39') '  -0.81%  '
Set-TextValue $ws.Range('D40') '1.204.67'
Set-TextValue $ws.Range('E40') '  -2.28%  '
Set-TextValue $ws.Range('D41') '6.422'
Set-TextValue $ws.Range('E41') '  -2.30%  '
Set-TextValue $ws.Range('D42') '0.9027'
Set-TextValue $ws.Range('E42') '  +0.74%  '
Set-TextValue $ws.Range('E43') '  +0.42%  '
Set-TextValue $ws.Range('D44') '101.36'
Set-TextValue $ws.Range('E44') '  -0.70%  '
Set-TextValue $ws.Range('D45') '1.989.02'
Set-TextValue $ws.Range('E45') '  +0.20%  '
Set-TextValue $ws.Range('E46') '  +0.40%  '
Set-TextValue $ws.Range('D47') '64.94'
Set-TextValue $ws.Range('E47') '  -1.39%  '
Set-TextValue $ws.Range('D48') '0.5134'
Set-TextValue $ws.Range('E48') '  +0.77%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D49') '0.4031'
Set-TextValue $ws.Range('E49') '  -0.83%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D50') '9.117'
Set-TextValue $ws.Range('E50') '  +0.67%  '
Set-TextValue $ws.Range('D51') '0.05833'
Set-TextValue $ws.Range('E51') '  +0.37%  '
